$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row of chart source data (row 2) is being updated: the doughnut
# chart now reports "Volunteering days by departments" instead of
# "Volunteering days by Directorates", with a new set of department
# labels replacing the old directorate labels.
$ws.Range("B2:J2").Select() | Out-Null

$ws.Range("A2").Value = "Volunteering days by departments for 2023-2024"
$ws.Range("B2").Value = "HR"
$ws.Range("C2").Value = "Finance"
$ws.Range("D2").Value = "IT"
$ws.Range("E2").Value = "Investment"
$ws.Range("F2").Value = "Customer Service"
$ws.Range("G2").Value = "Sale"
$ws.Range("H2").Value = "Legal"
$ws.Range("I2").Value = "Sustainability"
$ws.Range("J2").Value = "Diversity & Inclusion"

# AC2 (TOTAL) loses its explicit "General" number-format override and
# reverts to the sheet's default/no style.
$ws.Range("AC2").Style = "Normal"
